# Add a second table with per-month statistics (word count / entry count)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for the new table (row 68) ---
$ws.Range("F68").Value = "Month"

# --- Month rows (69-77) ---
$ws.Range("F69").Value = "March"
$ws.Range("F70").Value = "April"
$ws.Range("F71").Value = "May "
$ws.Range("F72").Value = "June"
$ws.Range("F73").Value = "July"
$ws.Range("F74").Value = "August"
$ws.Range("F75").Value = "September"
$ws.Range("F76").Value = "October"
$ws.Range("F77").Value = "November"

# --- Remaining header cells ---
$ws.Range("G68").Value = "# of words"
$ws.Range("H68").Value = "# of entries"

$ws.Range("G69").Formula = "=SUM(C2:C5)"
$ws.Range("H69").Value = 4

$ws.Range("G70").Formula = "=SUM(C6:C14)"
$ws.Range("H70").Value = 9

$ws.Range("G71").Formula = "=SUM(C15:C23)"
$ws.Range("H71").Value = 9

$ws.Range("G72").Formula = "=SUM(C24:C32)"
$ws.Range("H72").Value = 9

$ws.Range("G73").Formula = "=SUM(C33:C40)"
$ws.Range("H73").Value = 8

$ws.Range("G74").Formula = "=SUM(C41:C55)"
$ws.Range("H74").Value = 15

$ws.Range("G75").Formula = "=SUM(C56:C61)"
$ws.Range("H75").Value = 6

$ws.Range("G76").Formula = "=SUM(C62:C65)"
$ws.Range("H76").Value = 4

$ws.Range("G77").Formula = "=C66"
$ws.Range("H77").Value = 1

# --- Totals row (79) ---
$ws.Range("E79").Value = "Total"
$ws.Range("G79").Formula = "=SUM(G69:G77)"
$ws.Range("H79").Formula = "=SUM(H69:H77)"

# Move the active selection like the edited workbook
$ws.Range("E3").Select()
